# Day book: append 4 new transaction rows (row 3-6) below the existing header/data.
# Each row is first entered fully as text (date, time, id, work, input) - this mirrors
# how the source "day book" web app writes every form field as a string - and then the
# INPUT column is converted back to a real number for rows 3-5 (row 6 is left as text).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: 20/06/2018 06:47:09 SJ3 nbghb 52412
$ws.Range("A3").Value = "'20/06/2018"
$ws.Range("B3").Value = "'06:47:09"
$ws.Range("C3").Value = "'SJ3"
$ws.Range("D3").Value = "'nbghb"
$ws.Range("E3").Value = "'52412"

# Row 4: 20/06/2018 12:15:57 SJ3 asdsdgg 5000
$ws.Range("A4").Value = "'20/06/2018"
$ws.Range("B4").Value = "'12:15:57"
$ws.Range("C4").Value = "'SJ3"
$ws.Range("D4").Value = "'asdsdgg"
$ws.Range("E4").Value = "'5000"

# Row 5: 20/06/2018 12:22:59 SJ3 fdksdfbjm 100
$ws.Range("A5").Value = "'20/06/2018"
$ws.Range("B5").Value = "'12:22:59"
$ws.Range("C5").Value = "'SJ3"
$ws.Range("D5").Value = "'fdksdfbjm"
$ws.Range("E5").Value = "'100"

# Row 6: 20/06/2018 12:24:29 SJ3 nbnbv 600 (INPUT stays text for this row)
$ws.Range("A6").Value = "'20/06/2018"
$ws.Range("B6").Value = "'12:24:29"
$ws.Range("C6").Value = "'SJ3"
$ws.Range("D6").Value = "'nbnbv"
$ws.Range("E6").Value = "'600"

# Clear the forced-text formatting back to Normal on all the new cells.
$ws.Range("A3:E6").Style = "Normal"

# Convert the INPUT values of rows 3-5 into real numbers (row 6 keeps its text value).
$ws.Range("E3").Value = 52412
$ws.Range("E4").Value = 5000
$ws.Range("E5").Value = 100
